$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AD2").Value = 25
$ws.Range("AE2").Value = 2
$ws.Range("AI2").Value = 9
$ws.Range("AK2").Value = 17
$ws.Range("AN2").Value = 7
$ws.Range("AX2").Value = 5
$ws.Range("AY2").Value = 22
$ws.Range("AZ2").Value = 6
$ws.Range("BB2").Value = 18

$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 6
$ws.Range("G3").Value = 0.857
$ws.Range("I3").Value = 34.1
$ws.Range("J3").Value = 76.90000000000001
$ws.Range("L3").Value = 3.9
$ws.Range("N3").Value = 0.281
$ws.Range("O3").Value = 21.1
$ws.Range("P3").Value = 29.7
$ws.Range("Q3").Value = 0.712
$ws.Range("R3").Value = 11.7
$ws.Range("S3").Value = 34
$ws.Range("T3").Value = 45.7
$ws.Range("U3").Value = 19.9
$ws.Range("V3").Value = 18.6
$ws.Range("W3").Value = 9.6
$ws.Range("X3").Value = 5.4
$ws.Range("Y3").Value = 4.6
$ws.Range("Z3").Value = 25
$ws.Range("AA3").Value = 25.6
$ws.Range("AB3").Value = 93.3
$ws.Range("AC3").Value = 6.6
$ws.Range("AD3").Value = 1
$ws.Range("AK3").Value = 15
$ws.Range("AL3").Value = 28
$ws.Range("AN3").Value = 27
$ws.Range("AO3").Value = 9
$ws.Range("AQ3").Value = 24
$ws.Range("AR3").Value = 12
$ws.Range("AS3").Value = 4
$ws.Range("AT3").Value = 4
$ws.Range("AU3").Value = 17
$ws.Range("AX3").Value = 12
$ws.Range("AY3").Value = 11
$ws.Range("AZ3").Value = 27
$ws.Range("BB3").Value = 25

$ws.Range("AD4").Value = 11
$ws.Range("AE4").Value = 20
$ws.Range("AF4").Value = 18
$ws.Range("AG4").Value = 22
$ws.Range("AN4").Value = 15
$ws.Range("AO4").Value = 8
$ws.Range("AP4").Value = 11
$ws.Range("AQ4").Value = 8
$ws.Range("AS4").Value = 30
$ws.Range("AU4").Value = 27
$ws.Range("AV4").Value = 17
$ws.Range("AZ4").Value = 5
$ws.Range("BA4").Value = 18

$ws.Range("D5").Value = 7
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.429
$ws.Range("I5").Value = 33.9
$ws.Range("J5").Value = 80.59999999999999
$ws.Range("K5").Value = 0.42
$ws.Range("L5").Value = 5.6
$ws.Range("M5").Value = 14.6
$ws.Range("N5").Value = 0.382
$ws.Range("O5").Value = 22
$ws.Range("P5").Value = 28
$ws.Range("Q5").Value = 0.786
$ws.Range("R5").Value = 13
$ws.Range("S5").Value = 30
$ws.Range("T5").Value = 43
$ws.Range("U5").Value = 18.6
$ws.Range("V5").Value = 15.9
$ws.Range("W5").Value = 8.300000000000001
$ws.Range("Y5").Value = 5.1
$ws.Range("Z5").Value = 22.4
$ws.Range("AA5").Value = 23
$ws.Range("AB5").Value = 95.3
$ws.Range("AC5").Value = -0.3
$ws.Range("AD5").Value = 1
$ws.Range("AE5").Value = 13
$ws.Range("AF5").Value = 18
$ws.Range("AG5").Value = 16
$ws.Range("AI5").Value = 24
$ws.Range("AJ5").Value = 15
$ws.Range("AK5").Value = 26
$ws.Range("AL5").Value = 17
$ws.Range("AM5").Value = 22
$ws.Range("AO5").Value = 5
$ws.Range("AP5").Value = 9
$ws.Range("AQ5").Value = 9
$ws.Range("AS5").Value = 15
$ws.Range("AT5").Value = 10
$ws.Range("AU5").Value = 22
$ws.Range("AV5").Value = 21
$ws.Range("AW5").Value = 9
$ws.Range("AX5").Value = 4
$ws.Range("AY5").Value = 18
$ws.Range("AZ5").Value = 19
$ws.Range("BA5").Value = 11
$ws.Range("BB5").Value = 22
$ws.Range("BC5").Value = 16

$ws.Range("D6").Value = 7
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0.714
$ws.Range("I6").Value = 35.7
$ws.Range("J6").Value = 76.40000000000001
$ws.Range("K6").Value = 0.467
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 17.7
$ws.Range("N6").Value = 0.339
$ws.Range("O6").Value = 22.1
$ws.Range("P6").Value = 28.9
$ws.Range("Q6").Value = 0.767
$ws.Range("R6").Value = 10.9
$ws.Range("S6").Value = 30.1
$ws.Range("U6").Value = 19.4
$ws.Range("V6").Value = 15.3
$ws.Range("W6").Value = 7.9
$ws.Range("X6").Value = 6.1
$ws.Range("Y6").Value = 3.7
$ws.Range("Z6").Value = 22.1
$ws.Range("AA6").Value = 22.6
$ws.Range("AB6").Value = 99.59999999999999
$ws.Range("AC6").Value = 6.6
$ws.Range("AD6").Value = 1
$ws.Range("AE6").Value = 2
$ws.Range("AF6").Value = 5
$ws.Range("AG6").Value = 5
$ws.Range("AI6").Value = 17
$ws.Range("AJ6").Value = 27
$ws.Range("AK6").Value = 6
$ws.Range("AL6").Value = 14
$ws.Range("AM6").Value = 12
$ws.Range("AN6").Value = 18
$ws.Range("AP6").Value = 6
$ws.Range("AQ6").Value = 14
$ws.Range("AR6").Value = 21
$ws.Range("AS6").Value = 14
$ws.Range("AU6").Value = 20
$ws.Range("AV6").Value = 18
$ws.Range("AX6").Value = 7
$ws.Range("AY6").Value = 4
$ws.Range("AZ6").Value = 17
$ws.Range("BA6").Value = 13
$ws.Range("BB6").Value = 8
$ws.Range("BC6").Value = 5

$ws.Range("AD7").Value = 11
$ws.Range("AE7").Value = 20
$ws.Range("AF7").Value = 18
$ws.Range("AG7").Value = 22
$ws.Range("AI7").Value = 11
$ws.Range("AM7").Value = 9
$ws.Range("AO7").Value = 25
$ws.Range("AS7").Value = 9
$ws.Range("AT7").Value = 15
$ws.Range("AY7").Value = 8
$ws.Range("BA7").Value = 25
$ws.Range("BB7").Value = 20
$ws.Range("BC7").Value = 20

$ws.Range("AD8").Value = 11
$ws.Range("AF8").Value = 12
$ws.Range("AG8").Value = 13
$ws.Range("AH8").Value = 4
$ws.Range("AI8").Value = 18
$ws.Range("AJ8").Value = 9
$ws.Range("AP8").Value = 1
$ws.Range("AQ8").Value = 11
$ws.Range("AT8").Value = 3
$ws.Range("AV8").Value = 26
$ws.Range("AX8").Value = 10
$ws.Range("BC8").Value = 15

$ws.Range("AD9").Value = 11
$ws.Range("AE9").Value = 7
$ws.Range("AF9").Value = 5
$ws.Range("AI9").Value = 16
$ws.Range("AJ9").Value = 12
$ws.Range("AK9").Value = 16
$ws.Range("AL9").Value = 11
$ws.Range("AP9").Value = 16
$ws.Range("AT9").Value = 12
$ws.Range("AU9").Value = 9
$ws.Range("AX9").Value = 9
$ws.Range("AY9").Value = 6
$ws.Range("AZ9").Value = 20
$ws.Range("BA9").Value = 14

$ws.Range("AE10").Value = 20
$ws.Range("AO10").Value = 6
$ws.Range("AU10").Value = 23
$ws.Range("AV10").Value = 13
$ws.Range("BB10").Value = 5
$ws.Range("BC10").Value = 21

$ws.Range("AE11").Value = 7
$ws.Range("AF11").Value = 12
$ws.Range("AG11").Value = 12
$ws.Range("AH11").Value = 7
$ws.Range("AJ11").Value = 19
$ws.Range("AK11").Value = 27
$ws.Range("AL11").Value = 13
$ws.Range("AM11").Value = 13
$ws.Range("AN11").Value = 13
$ws.Range("AO11").Value = 6
$ws.Range("AP11").Value = 15
$ws.Range("AR11").Value = 19
$ws.Range("AS11").Value = 11
$ws.Range("AT11").Value = 13
$ws.Range("AU11").Value = 26
$ws.Range("BA11").Value = 15

$ws.Range("AD12").Value = 25
$ws.Range("AE12").Value = 20
$ws.Range("AF12").Value = 12
$ws.Range("AG12").Value = 20
$ws.Range("AI12").Value = 13
$ws.Range("AJ12").Value = 14
$ws.Range("AL12").Value = 4
$ws.Range("AN12").Value = 16
$ws.Range("AR12").Value = 11
$ws.Range("AS12").Value = 8
$ws.Range("AT12").Value = 9
$ws.Range("AU12").Value = 10
$ws.Range("AX12").Value = 11
$ws.Range("BA12").Value = 15

$ws.Range("AH13").Value = 7
$ws.Range("AJ13").Value = 10
$ws.Range("AK13").Value = 28
$ws.Range("AM13").Value = 14
$ws.Range("AP13").Value = 26
$ws.Range("AR13").Value = 18
$ws.Range("AS13").Value = 15
$ws.Range("AU13").Value = 23
$ws.Range("AW13").Value = 22
$ws.Range("AX13").Value = 7
$ws.Range("AZ13").Value = 26
$ws.Range("BA13").Value = 24

$ws.Range("D14").Value = 5
$ws.Range("E14").Value = 5
$ws.Range("I14").Value = 37.2
$ws.Range("J14").Value = 82
$ws.Range("K14").Value = 0.454
$ws.Range("L14").Value = 8
$ws.Range("M14").Value = 17
$ws.Range("N14").Value = 0.471
$ws.Range("O14").Value = 24.4
$ws.Range("P14").Value = 32.4
$ws.Range("Q14").Value = 0.753
$ws.Range("R14").Value = 13.4
$ws.Range("S14").Value = 38
$ws.Range("T14").Value = 51.4
$ws.Range("U14").Value = 21
$ws.Range("V14").Value = 15.4
$ws.Range("W14").Value = 9.6
$ws.Range("Y14").Value = 4.4
$ws.Range("Z14").Value = 19.4
$ws.Range("AA14").Value = 23.8
$ws.Range("AB14").Value = 106.8
$ws.Range("AC14").Value = 22.4
$ws.Range("AD14").Value = 25
$ws.Range("AE14").Value = 2
$ws.Range("AK14").Value = 9
$ws.Range("AM14").Value = 15
$ws.Range("AN14").Value = 1
$ws.Range("AP14").Value = 2
$ws.Range("AQ14").Value = 19
$ws.Range("AR14").Value = 4
$ws.Range("AU14").Value = 10
$ws.Range("AV14").Value = 19
$ws.Range("AW14").Value = 3
$ws.Range("AX14").Value = 3
$ws.Range("AY14").Value = 10
$ws.Range("AZ14").Value = 4
$ws.Range("BA14").Value = 4

$ws.Range("AF15").Value = 18
$ws.Range("AG15").Value = 16
$ws.Range("AJ15").Value = 24
$ws.Range("AP15").Value = 12
$ws.Range("AR15").Value = 17
$ws.Range("AX15").Value = 15
$ws.Range("BA15").Value = 17
$ws.Range("BC15").Value = 19

$ws.Range("AD16").Value = 11
$ws.Range("AF16").Value = 12
$ws.Range("AG16").Value = 13
$ws.Range("AI16").Value = 15
$ws.Range("AK16").Value = 8
$ws.Range("AM16").Value = 8
$ws.Range("AN16").Value = 17
$ws.Range("AO16").Value = 11
$ws.Range("AX16").Value = 14
$ws.Range("BB16").Value = 7

$ws.Range("D17").Value = 7
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 0.429
$ws.Range("H17").Value = 48.7
$ws.Range("J17").Value = 81.59999999999999
$ws.Range("K17").Value = 0.45
$ws.Range("L17").Value = 5.1
$ws.Range("N17").Value = 0.375
$ws.Range("O17").Value = 17.3
$ws.Range("P17").Value = 24.4
$ws.Range("Q17").Value = 0.708
$ws.Range("R17").Value = 11.6
$ws.Range("S17").Value = 29.9
$ws.Range("T17").Value = 41.4
$ws.Range("U17").Value = 22.7
$ws.Range("V17").Value = 17.3
$ws.Range("W17").Value = 6.6
$ws.Range("X17").Value = 3.6
$ws.Range("Y17").Value = 4.7
$ws.Range("Z17").Value = 26.4
$ws.Range("AA17").Value = 23.3
$ws.Range("AB17").Value = 95.90000000000001
$ws.Range("AC17").Value = -1.4
$ws.Range("AD17").Value = 1
$ws.Range("AF17").Value = 18
$ws.Range("AG17").Value = 16
$ws.Range("AH17").Value = 7
$ws.Range("AI17").Value = 10
$ws.Range("AK17").Value = 11
$ws.Range("AN17").Value = 11
$ws.Range("AP17").Value = 19
$ws.Range("AQ17").Value = 25
$ws.Range("AR17").Value = 14
$ws.Range("AS17").Value = 18
$ws.Range("AT17").Value = 16
$ws.Range("AV17").Value = 27
$ws.Range("AW17").Value = 22
$ws.Range("AX17").Value = 26
$ws.Range("AY17").Value = 13
$ws.Range("BA17").Value = 9
$ws.Range("BC17").Value = 17

$ws.Range("AD18").Value = 11
$ws.Range("AM18").Value = 21
$ws.Range("AP18").Value = 17
$ws.Range("AR18").Value = 15
$ws.Range("AY18").Value = 20
$ws.Range("AZ18").Value = 18
$ws.Range("BC18").Value = 22

$ws.Range("AD19").Value = 25
$ws.Range("AE19").Value = 20
$ws.Range("AF19").Value = 12
$ws.Range("AG19").Value = 20
$ws.Range("AJ19").Value = 18
$ws.Range("AK19").Value = 29
$ws.Range("AL19").Value = 16
$ws.Range("AM19").Value = 11
$ws.Range("AO19").Value = 10
$ws.Range("AP19").Value = 13
$ws.Range("AQ19").Value = 10
$ws.Range("AR19").Value = 16
$ws.Range("AU19").Value = 28
$ws.Range("AX19").Value = 13

$ws.Range("AD20").Value = 11
$ws.Range("AE20").Value = 7
$ws.Range("AF20").Value = 5
$ws.Range("AJ20").Value = 26
$ws.Range("AP20").Value = 17
$ws.Range("AQ20").Value = 13

$ws.Range("AD21").Value = 11
$ws.Range("AE21").Value = 7
$ws.Range("AF21").Value = 5
$ws.Range("AN21").Value = 12
$ws.Range("AO21").Value = 23
$ws.Range("AQ21").Value = 20
$ws.Range("AR21").Value = 19
$ws.Range("AS21").Value = 19
$ws.Range("AV21").Value = 15
$ws.Range("AW21").Value = 2
$ws.Range("AY21").Value = 6

$ws.Range("AD22").Value = 11
$ws.Range("AI22").Value = 25
$ws.Range("AN22").Value = 14
$ws.Range("AO22").Value = 26
$ws.Range("AP22").Value = 25
$ws.Range("AQ22").Value = 23
$ws.Range("AS22").Value = 15
$ws.Range("AT22").Value = 14
$ws.Range("AV22").Value = 20
$ws.Range("AW22").Value = 13
$ws.Range("AY22").Value = 16
$ws.Range("AZ22").Value = 21

$ws.Range("AD23").Value = 11
$ws.Range("AE23").Value = 7
$ws.Range("AF23").Value = 5
$ws.Range("AI23").Value = 18
$ws.Range("AL23").Value = 4
$ws.Range("AP23").Value = 8
$ws.Range("AR23").Value = 7
$ws.Range("AS23").Value = 9
$ws.Range("AU23").Value = 25
$ws.Range("AV23").Value = 14
$ws.Range("AY23").Value = 8
$ws.Range("BA23").Value = 8
$ws.Range("BC23").Value = 4

$ws.Range("AD24").Value = 11
$ws.Range("AE24").Value = 20
$ws.Range("AF24").Value = 18
$ws.Range("AG24").Value = 22
$ws.Range("AO24").Value = 23
$ws.Range("AQ24").Value = 16
$ws.Range("AS24").Value = 3
$ws.Range("AY24").Value = 15
$ws.Range("BB24").Value = 17

$ws.Range("AE25").Value = 2
$ws.Range("AF25").Value = 5
$ws.Range("AG25").Value = 5
$ws.Range("AI25").Value = 12
$ws.Range("AQ25").Value = 15
$ws.Range("AY25").Value = 5
$ws.Range("BA25").Value = 7

$ws.Range("D26").Value = 6
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 0.5
$ws.Range("H26").Value = 48.8
$ws.Range("I26").Value = 35.2
$ws.Range("J26").Value = 80.3
$ws.Range("K26").Value = 0.438
$ws.Range("L26").Value = 8
$ws.Range("M26").Value = 20.2
$ws.Range("N26").Value = 0.397
$ws.Range("O26").Value = 16
$ws.Range("P26").Value = 21
$ws.Range("Q26").Value = 0.762
$ws.Range("R26").Value = 13.7
$ws.Range("S26").Value = 25.8
$ws.Range("T26").Value = 39.5
$ws.Range("U26").Value = 19.8
$ws.Range("V26").Value = 11.7
$ws.Range("W26").Value = 6.8
$ws.Range("X26").Value = 5
$ws.Range("Y26").Value = 2.7
$ws.Range("Z26").Value = 21.7
$ws.Range("AA26").Value = 19.3
$ws.Range("AB26").Value = 94.3
$ws.Range("AC26").Value = -5.2
$ws.Range("AD26").Value = 11
$ws.Range("AF26").Value = 12
$ws.Range("AG26").Value = 13
$ws.Range("AH26").Value = 4
$ws.Range("AI26").Value = 20
$ws.Range("AJ26").Value = 17
$ws.Range("AK26").Value = 18
$ws.Range("AL26").Value = 4
$ws.Range("AM26").Value = 6
$ws.Range("AN26").Value = 6
$ws.Range("AO26").Value = 26
$ws.Range("AP26").Value = 27
$ws.Range("AQ26").Value = 17
$ws.Range("AR26").Value = 3
$ws.Range("AS26").Value = 29
$ws.Range("AT26").Value = 23
$ws.Range("AU26").Value = 18
$ws.Range("AW26").Value = 19
$ws.Range("AX26").Value = 15
$ws.Range("AZ26").Value = 15
$ws.Range("BA26").Value = 25
$ws.Range("BB26").Value = 23
$ws.Range("BC26").Value = 22

$ws.Range("AF27").Value = 18
$ws.Range("AG27").Value = 16
$ws.Range("AU27").Value = 16
$ws.Range("AV27").Value = 25
$ws.Range("AY27").Value = 13
$ws.Range("BB27").Value = 6

$ws.Range("AD28").Value = 25
$ws.Range("AF28").Value = 18
$ws.Range("AK28").Value = 4
$ws.Range("AL28").Value = 8
$ws.Range("AM28").Value = 7
$ws.Range("AS28").Value = 13
$ws.Range("AT28").Value = 24
$ws.Range("AU28").Value = 15
$ws.Range("AY28").Value = 12

$ws.Range("AD29").Value = 11
$ws.Range("AE29").Value = 7
$ws.Range("AF29").Value = 5
$ws.Range("AH29").Value = 4
$ws.Range("AK29").Value = 5
$ws.Range("AL29").Value = 11
$ws.Range("AN29").Value = 2
$ws.Range("AS29").Value = 19
$ws.Range("AX29").Value = 6
$ws.Range("AY29").Value = 20
$ws.Range("BB29").Value = 21

$ws.Range("AD30").Value = 11
$ws.Range("AE30").Value = 2
$ws.Range("AG30").Value = 4
$ws.Range("AN30").Value = 26
$ws.Range("AO30").Value = 13
$ws.Range("AP30").Value = 14
$ws.Range("AQ30").Value = 12
$ws.Range("AS30").Value = 12
$ws.Range("AT30").Value = 10
$ws.Range("AV30").Value = 15
$ws.Range("AY30").Value = 16
$ws.Range("BA30").Value = 6
$ws.Range("BC30").Value = 3

$ws.Range("AD31").Value = 25
$ws.Range("AI31").Value = 13
$ws.Range("AJ31").Value = 13
$ws.Range("AO31").Value = 12
$ws.Range("AP31").Value = 7
$ws.Range("AQ31").Value = 26
$ws.Range("AT31").Value = 24
$ws.Range("AU31").Value = 19
$ws.Range("AV31").Value = 12
$ws.Range("BA31").Value = 4
$ws.Range("BB31").Value = 16
